$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain literal text
# (matches the source file, which stores every Price cell as inline text)
$textCells = @("D5", "D7", "D10", "D13", "D14", "D16", "D20", "D22", "D26", "D27", "D30", "D33", "D35", "D36", "D37", "D38", "D39", "D42", "D43", "D46", "D48", "D49")
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row
$ws.Range("D2").Value = "38.149.61"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.092.47"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "228.80"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "60.76"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "0.0849"
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "2.402.83"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "14.66"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "22.22"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("D16").Value = "0.777"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "2.091.75"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "38.066.05"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "70.13"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "223.88"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "169.86"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +6.40%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("D33").Value = "4.70"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "0.0607"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").Value = "6.38"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "3.51"
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").Value = "1.555.37"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "100.10"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").Value = "0.0219"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "7.46"
$ws.Range("E48").Value = "  +4.93%  "
$ws.Range("D49").Value = "1.03"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "2.289.12"
$ws.Range("E51").Value = "  +2.71%  "
